$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 4759.8
$ws.Range("I32").Value2 = 4500
$ws.Range("J32").Value2 = 4824.75
$ws.Range("K32").Value2 = 4500
$ws.Range("L32").Value2 = 4824.75
$ws.Range("M32").Value2 = -4174
$ws.Range("N32").Value2 = -5476.75

$ws.Range("H51").Value2 = 3800.3333
$ws.Range("I51").Value2 = 2200.5
$ws.Range("K51").Value2 = 2200.5
$ws.Range("M51").Value2 = -1716.5

$ws.Range("H62").Value2 = 14139.9
$ws.Range("I62").Value2 = 11628.429
$ws.Range("K62").Value2 = 11628.429
$ws.Range("M62").Value2 = -11004.429

$ws.Range("H64").Value2 = 8818.6
$ws.Range("I64").Value2 = 4999.5
$ws.Range("J64").Value2 = 9242.944
$ws.Range("K64").Value2 = 4999.5
$ws.Range("L64").Value2 = 9242.944
$ws.Range("M64").Value2 = -4751.5
$ws.Range("N64").Value2 = -9738.944

$ws.Range("H65").Value2 = 14139.9
$ws.Range("I65").Value2 = 11628.429
$ws.Range("K65").Value2 = 58142.145
$ws.Range("M65").Value2 = -55022.145

$ws.Range("H67").Value2 = 8818.6
$ws.Range("I67").Value2 = 4999.5
$ws.Range("J67").Value2 = 9242.944
$ws.Range("K67").Value2 = 4999.5
$ws.Range("L67").Value2 = 9242.944
$ws.Range("M67").Value2 = -4141.5
$ws.Range("N67").Value2 = -10958.944

$ws.Range("H70").Value2 = 4013.2856
$ws.Range("I70").Value2 = 1780
$ws.Range("J70").Value2 = 5254
$ws.Range("K70").Value2 = 5340
$ws.Range("L70").Value2 = 15762
$ws.Range("M70").Value2 = -5070
$ws.Range("N70").Value2 = -16302

$ws.Range("H73").Value2 = 4013.2856
$ws.Range("I73").Value2 = 1780
$ws.Range("J73").Value2 = 5254
$ws.Range("K73").Value2 = 5340
$ws.Range("L73").Value2 = 15762
$ws.Range("M73").Value2 = -4404
$ws.Range("N73").Value2 = -17634

$ws.Range("H74").Value2 = 6891.8
$ws.Range("I74").Value2 = 4093.7778
$ws.Range("K74").Value2 = 4093.7778
$ws.Range("M74").Value2 = -3157.7778

$ws.Range("H77").Value2 = 6891.8
$ws.Range("I77").Value2 = 4093.7778
$ws.Range("K77").Value2 = 20468.889
$ws.Range("M77").Value2 = -15788.889

$ws.Range("H82").Value2 = 4677.5557
$ws.Range("I82").Value2 = 4677.5557
$ws.Range("K82").Value2 = 14032.6671
$ws.Range("M82").Value2 = -13626.6671

$ws.Range("H85").Value2 = 4677.5557
$ws.Range("I85").Value2 = 4677.5557
$ws.Range("K85").Value2 = 14032.6671
$ws.Range("M85").Value2 = -12628.6671

$ws.Range("H86").Value2 = 2097.5
$ws.Range("I86").Value2 = 1821.1
$ws.Range("K86").Value2 = 1821.1
$ws.Range("M86").Value2 = -698.0999999999999

$ws.Range("H87").Value2 = 97994.60000000001
$ws.Range("I87").Value2 = 90000
$ws.Range("J87").Value2 = 99993.25
$ws.Range("K87").Value2 = 90000
$ws.Range("L87").Value2 = 99993.25
$ws.Range("M87").Value2 = -88752
$ws.Range("N87").Value2 = -102489.25

$ws.Range("H88").Value2 = 22060
$ws.Range("J88").Value2 = 24986.666
$ws.Range("L88").Value2 = 24986.666
$ws.Range("N88").Value2 = -25798.666

$ws.Range("H89").Value2 = 2097.5
$ws.Range("I89").Value2 = 1821.1
$ws.Range("K89").Value2 = 9105.5
$ws.Range("M89").Value2 = -3489.5

$ws.Range("H90").Value2 = 97994.60000000001
$ws.Range("I90").Value2 = 90000
$ws.Range("J90").Value2 = 99993.25
$ws.Range("K90").Value2 = 270000
$ws.Range("L90").Value2 = 299979.75
$ws.Range("M90").Value2 = -263760
$ws.Range("N90").Value2 = -312459.75

$ws.Range("H91").Value2 = 22060
$ws.Range("J91").Value2 = 24986.666
$ws.Range("L91").Value2 = 24986.666
$ws.Range("N91").Value2 = -27794.666

$ws.Range("H92").Value2 = 1156.8182
$ws.Range("I92").Value2 = 272
$ws.Range("K92").Value2 = 272
$ws.Range("M92").Value2 = 976

$ws.Range("H96").Value2 = 2187.1428
$ws.Range("I96").Value2 = 958.5454999999999
$ws.Range("K96").Value2 = 2875.6365
$ws.Range("M96").Value2 = -1502.6365

$ws.Range("H101").Value2 = 1928.4286
$ws.Range("I101").Value2 = 1416.5
$ws.Range("K101").Value2 = 4249.5
$ws.Range("M101").Value2 = -2627.5

$ws.Range("H116").Value2 = 268106.2
$ws.Range("I116").Value2 = 123250
$ws.Range("J116").Value2 = 333949.9
$ws.Range("K116").Value2 = 123250
$ws.Range("L116").Value2 = 333949.9
$ws.Range("M116").Value2 = -119808
$ws.Range("N116").Value2 = -340833.9

$ws.Range("H120").Value2 = 88999
$ws.Range("J120").Value2 = 88999
$ws.Range("L120").Value2 = 88999
$ws.Range("N120").Value2 = -98675

$ws.Range("H135").Value2 = 2137.7646
$ws.Range("I135").Value2 = 680.38464
$ws.Range("K135").Value2 = 6123.46176
$ws.Range("M135").Value2 = -3588.46176

$ws.Range("H137").Value2 = 2957739.8
$ws.Range("I137").Value2 = 7865.3335
$ws.Range("K137").Value2 = 23596.0005
$ws.Range("M137").Value2 = -21046.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 11100.765
$ws.Range("I2").Value2 = 2276.25
$ws.Range("K2").Value2 = 2276.25
$ws.Range("M2").Value2 = -2163.25

$ws.Range("H45").Value2 = 1849.8235
$ws.Range("I45").Value2 = 1860.5714
$ws.Range("J45").Value2 = 1799.6666
$ws.Range("K45").Value2 = 1860.5714
$ws.Range("L45").Value2 = 1799.6666
$ws.Range("M45").Value2 = -1483.5714
$ws.Range("N45").Value2 = -2553.6666

$ws.Range("H61").Value2 = 1591407.6
$ws.Range("I61").Value2 = 2384328.5
$ws.Range("K61").Value2 = 2384328.5
$ws.Range("M61").Value2 = -2384116.5

$ws.Range("H116").Value2 = 11100.765
$ws.Range("I116").Value2 = 2276.25
$ws.Range("K116").Value2 = 2276.25
$ws.Range("M116").Value2 = 17.75

$ws.Range("H132").Value2 = 610691.7
$ws.Range("I132").Value2 = 743174.4
$ws.Range("K132").Value2 = 2229523.2
$ws.Range("M132").Value2 = -2226993.2

$ws.Range("H136").Value2 = 1591407.6
$ws.Range("I136").Value2 = 2384328.5
$ws.Range("K136").Value2 = 7152985.5
$ws.Range("M136").Value2 = -7150435.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 11100.765
$ws.Range("I3").Value2 = 2276.25
$ws.Range("K3").Value2 = 2276.25
$ws.Range("M3").Value2 = -2162.25

$ws.Range("H82").Value2 = 67316.875
$ws.Range("I82").Value2 = 31628.5
$ws.Range("K82").Value2 = 31628.5
$ws.Range("M82").Value2 = -31245.5

$ws.Range("H85").Value2 = 67316.875
$ws.Range("I85").Value2 = 31628.5
$ws.Range("K85").Value2 = 31628.5
$ws.Range("M85").Value2 = -30302.5

$ws.Range("H107").Value2 = 9491.111000000001
$ws.Range("I107").Value2 = 10177.625
$ws.Range("K107").Value2 = 10177.625
$ws.Range("M107").Value2 = -8257.625

$ws.Range("H134").Value2 = 1695007
$ws.Range("I134").Value2 = 6041167
$ws.Range("J134").Value2 = 357726.94
$ws.Range("K134").Value2 = 18123501
$ws.Range("L134").Value2 = 1073180.82
$ws.Range("M134").Value2 = -18120966
$ws.Range("N134").Value2 = -1078250.82

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 119051.7
$ws.Range("I31").Value2 = 241728.84
$ws.Range("J31").Value2 = 25239.766
$ws.Range("K31").Value2 = 241728.84
$ws.Range("L31").Value2 = 25239.766
$ws.Range("M31").Value2 = -241433.84
$ws.Range("N31").Value2 = -25829.766

$ws.Range("H34").Value2 = 119051.7
$ws.Range("I34").Value2 = 241728.84
$ws.Range("J34").Value2 = 25239.766
$ws.Range("K34").Value2 = 241728.84
$ws.Range("L34").Value2 = 25239.766
$ws.Range("M34").Value2 = -241526.84
$ws.Range("N34").Value2 = -25643.766

$ws.Range("H86").Value2 = 5821.125
$ws.Range("I86").Value2 = 5862.222
$ws.Range("J86").Value2 = 5768.2856
$ws.Range("K86").Value2 = 5862.222
$ws.Range("L86").Value2 = 5768.2856
$ws.Range("M86").Value2 = -4739.222
$ws.Range("N86").Value2 = -8014.2856

$ws.Range("H89").Value2 = 5821.125
$ws.Range("I89").Value2 = 5862.222
$ws.Range("J89").Value2 = 5768.2856
$ws.Range("K89").Value2 = 29311.11
$ws.Range("L89").Value2 = 28841.428
$ws.Range("M89").Value2 = -23695.11
$ws.Range("N89").Value2 = -40073.428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value2 = 44999
$ws.Range("J52").Value2 = 44999
$ws.Range("L52").Value2 = 44999
$ws.Range("N52").Value2 = -45517

$ws.Range("H102").Value2 = 3622.2068
$ws.Range("I102").Value2 = 2814.4167
$ws.Range("K102").Value2 = 2814.4167
$ws.Range("M102").Value2 = -1192.4167

$ws.Range("H122").Value2 = 56174
$ws.Range("I122").Value2 = 103263.2
$ws.Range("J122").Value2 = 9084.799999999999
$ws.Range("K122").Value2 = 309789.6
$ws.Range("L122").Value2 = 27254.4
$ws.Range("M122").Value2 = -307339.6
$ws.Range("N122").Value2 = -32154.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value2 = 500000160
$ws.Range("J2").Value2 = 500000160
$ws.Range("L2").Value2 = 500000160
$ws.Range("N2").Value2 = -500000384

$ws.Range("H4").Value2 = 91141200
$ws.Range("J4").Value2 = 100005320
$ws.Range("L4").Value2 = 100005320
$ws.Range("N4").Value2 = -100005546

$ws.Range("H6").Value2 = 640.2
$ws.Range("J6").Value2 = 687.75
$ws.Range("L6").Value2 = 687.75
$ws.Range("N6").Value2 = -917.75

$ws.Range("H8").Value2 = 17583
$ws.Range("I8").Value2 = 12750
$ws.Range("J8").Value2 = 19999.5
$ws.Range("K8").Value2 = 12750
$ws.Range("L8").Value2 = 19999.5
$ws.Range("M8").Value2 = -12610
$ws.Range("N8").Value2 = -20279.5

$ws.Range("H62").Value2 = 153199.83
$ws.Range("I62").Value2 = 451600
$ws.Range("J62").Value2 = 3999.75
$ws.Range("K62").Value2 = 451600
$ws.Range("L62").Value2 = 3999.75
$ws.Range("M62").Value2 = -450976
$ws.Range("N62").Value2 = -5247.75

$ws.Range("H65").Value2 = 153199.83
$ws.Range("I65").Value2 = 451600
$ws.Range("J65").Value2 = 3999.75
$ws.Range("K65").Value2 = 2258000
$ws.Range("L65").Value2 = 19998.75
$ws.Range("M65").Value2 = -2254880
$ws.Range("N65").Value2 = -26238.75

$ws.Range("H96").Value2 = 1958.3334
$ws.Range("J96").Value2 = 2000
$ws.Range("L96").Value2 = 2000
$ws.Range("N96").Value2 = -4746
